# Refresh the cryptocurrency snapshot (prices / 1h volume %) and fix
# two rows whose coin name + link had been swapped, per the latest
# coinranking.com pull performed by the scheduled GitHub Action.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.066.25"
$ws.Range("E2").Value = "  -0.67%  "

$ws.Range("D3").Value = "2.238.85"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("E4").Value = "  +0.01%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "242.96"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.90%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.624"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.73%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "74.22"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("E8").Value = "  +0.07%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.598"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -2.62%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "42.09"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -2.71%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0952"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -0.66%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.103"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +0.10%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "6.92"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -2.76%  "

$ws.Range("D14").Value = "2.573.88"
$ws.Range("E14").Value = "  +0.13%  "

$ws.Range("E15").Value = "  -0.24%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.838"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -1.69%  "

$ws.Range("D17").Value = "2.237.48"
$ws.Range("E17").Value = "  +1.45%  "

$ws.Range("D18").Value = "41.984.26"
$ws.Range("E18").Value = "  -0.60%  "

$ws.Range("E19").Value = "  -2.69%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "6.22"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.70%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "72.64"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.76%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "11.25"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +9.46%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "229.88"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.73%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.04"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -6.45%  "

$ws.Range("E25").Value = "  +0.00%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "11.39"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -3.30%  "

$ws.Range("E27").Value = "  -0.60%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.27"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -0.89%  "

$ws.Range("E29").Value = "  -1.16%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "167.79"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +0.62%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "20.58"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -1.62%  "

$ws.Range("E32").Value = "  -4.36%  "

$ws.Range("E33").Value = "  -0.72%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "30.06"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +0.88%  "

$ws.Range("E35").Value = "  -0.23%  "

$ws.Range("E36").Value = "  -6.52%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "4.28"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -4.37%  "

$ws.Range("E38").Value = "  -1.73%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "13.13"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -0.78%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.13"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -1.69%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "5.69"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.82%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "64.77"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +1.96%  "

$ws.Range("E43").Value = "  -1.35%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "8.71"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -2.12%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "103.78"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -2.22%  "

$ws.Range("E46").Value = "  -1.92%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "1.13"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -0.35%  "

$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.17"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.86%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "2.33"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -2.24%  "

$ws.Range("E50").Value = "  -2.15%  "

$ws.Range("D51").Value = "2.448.04"
$ws.Range("E51").Value = "  -0.03%  "
